$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date-format style from BA1 into BB1, then set the new period-end date value
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# New column BB ("latest" quarter) duplicates column BA for existing rows 2-81
# (EQUIPMENT eval carried forward), with updated figures for the most recent quarters.
$ws.Range("BB2").Value = 2.129035685702092
$ws.Range("BB3").Value = 2.893366588103589
$ws.Range("BB4").Value = 2.955627504208366
$ws.Range("BB5").Value = 1.167401765354654
$ws.Range("BB6").Value = 4.206619947727546
$ws.Range("BB7").Value = 2.391456992346349
$ws.Range("BB8").Value = 2.759646870517571
$ws.Range("BB9").Value = 6.320819196802205
$ws.Range("BB10").Value = -1.962345281877049
$ws.Range("BB11").Value = 2.026397378686767
$ws.Range("BB12").Value = 1.814496316115338
$ws.Range("BB13").Value = 1.63165117816655
$ws.Range("BB14").Value = 1.358419700277324
$ws.Range("BB15").Value = -0.1398451755242718
$ws.Range("BB16").Value = -0.07003068004286206
$ws.Range("BB17").Value = -7.3868590581191
$ws.Range("BB18").Value = -10.33903131837766
$ws.Range("BB19").Value = -1.060235412222937
$ws.Range("BB20").Value = 3.423433284297019
$ws.Range("BB21").Value = 2.063081011733999
$ws.Range("BB22").Value = 2.99202665168275
$ws.Range("BB23").Value = 6.534207423404695
$ws.Range("BB24").Value = 1.294218807309846
$ws.Range("BB25").Value = 2.007778863461724
$ws.Range("BB26").Value = 2.488825190489734
$ws.Range("BB27").Value = -0.02236378853282872
$ws.Range("BB28").Value = 1.864148130528193
$ws.Range("BB29").Value = 0.3074416423962276
$ws.Range("BB30").Value = 1.591452302439862
$ws.Range("BB31").Value = 1.24229439238384
$ws.Range("BB32").Value = 1.021051258256691
$ws.Range("BB33").Value = -2.083516682089652
$ws.Range("BB34").Value = 0.8945075486844729
$ws.Range("BB35").Value = 1.373553924371535
$ws.Range("BB36").Value = 0.3276711086097635
$ws.Range("BB37").Value = 1.809509533486136
$ws.Range("BB38").Value = 0.4339813219297071
$ws.Range("BB39").Value = 0.4041977135476031
$ws.Range("BB40").Value = 1.507492882068462
$ws.Range("BB41").Value = 1.350099422002103
$ws.Range("BB42").Value = 1.114800352984872
$ws.Range("BB43").Value = 1.516137977472326
$ws.Range("BB44").Value = -0.02409266101658147
$ws.Range("BB45").Value = -0.7388860519741201
$ws.Range("BB46").Value = 1.567695497950282
$ws.Range("BB47").Value = 1.091382109524247
$ws.Range("BB48").Value = -0.2679292251141305
$ws.Range("BB49").Value = 1.682996656587392
$ws.Range("BB50").Value = 1.591260107191601
$ws.Range("BB51").Value = 1.13706115148689
$ws.Range("BB52").Value = 1.735453665039003
$ws.Range("BB53").Value = 2.625533283765208
$ws.Range("BB54").Value = -0.3243298885145123
$ws.Range("BB55").Value = 0.6372341859553217
$ws.Range("BB56").Value = -0.5789332341234967
$ws.Range("BB57").Value = 0.1728728569232914
$ws.Range("BB58").Value = 1.7
$ws.Range("BB59").Value = -1.4
$ws.Range("BB60").Value = 1.3
$ws.Range("BB61").Value = -0.6
$ws.Range("BB62").Value = -3.292009884772611
$ws.Range("BB63").Value = -20.353166912592
$ws.Range("BB64").Value = 17.87631681612835
$ws.Range("BB65").Value = 4.395171409529297
$ws.Range("BB66").Value = 1.763269038133103
$ws.Range("BB67").Value = 1.651360288740889
$ws.Range("BB68").Value = -0.6023499639201475
$ws.Range("BB69").Value = 3.773751128807561
$ws.Range("BB70").Value = 0.3486937127494798
$ws.Range("BB71").Value = 0.6021692060546258
$ws.Range("BB72").Value = 1.505090528591751
$ws.Range("BB73").Value = -1.274647749701984
$ws.Range("BB74").Value = 0.4647476577308112
$ws.Range("BB75").Value = -0.6885297541760451
$ws.Range("BB76").Value = -0.9272934801906558
$ws.Range("BB77").Value = -0.8572418502401149
$ws.Range("BB78").Value = 0.05364604092510206
$ws.Range("BB79").Value = 1.830682919206694
$ws.Range("BB80").Value = -2.604576060482884
$ws.Range("BB81").Value = -3.124437332092583

# Rows 82-83: refreshed values for column BB (revised vs BA)
$ws.Range("BB82").Value = 2.411315004676197
$ws.Range("BB83").Value = 0.2886831937783967

# New row 84: next reporting period date (copy date style from A83) and BB value
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A84").Value = 45884
$ws.Range("BB84").Value = -0.7196185376451893
